# Results_tdd/tests_timings.xlsx is a small manual test-timings sheet. The
# commit ("add python script for hypothesis test") mainly adds a script
# elsewhere in the repo; the only change to this workbook is that the
# now-unused "average" summary row (row 9 — which held =AVERAGE(...)
# formulas across all 8 data columns) was cleared out before saving, and
# the selection/cursor was left parked on B12.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the =AVERAGE(...) formulas (and their cached results) from row 9,
# columns A:H, while leaving the existing cell formatting in place.
$ws.Range("A9:H9").ClearContents()

# Park the selection on B12, matching the saved workbook state.
$ws.Range("B12").Select()
